$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry below is a cell reference paired with its new text value, taken
# from the latest cryptos snapshot.
#
# Column D holds price figures that are stored as TEXT in the source sheet
# (e.g. "63.870.53", "2.00"), not as numbers. When such a purely-numeric-
# looking string (like "566.45") is assigned through .Value, Excel would
# normally reinterpret it as a real number and lose the original text
# formatting (trailing zeros, etc.). To avoid that we prefix those values
# with a leading apostrophe, exactly as typing '566.45 into a cell does,
# which keeps Excel from converting it to a number.
$updates = @(
    @{ Cell = 'D2'; Value = '63.870.53' }
    @{ Cell = 'E2'; Value = '  +0.04%  ' }
    @{ Cell = 'D3'; Value = '2.737.93' }
    @{ Cell = 'E3'; Value = '  -0.50%  ' }
    @{ Cell = 'E4'; Value = '  +0.02%  ' }
    @{ Cell = 'D5'; Value = '566.45' }
    @{ Cell = 'E5'; Value = '  -1.05%  ' }
    @{ Cell = 'D6'; Value = '161.35' }
    @{ Cell = 'E6'; Value = '  +2.80%  ' }
    @{ Cell = 'D7'; Value = '0.999' }
    @{ Cell = 'E7'; Value = '  +0.00%  ' }
    @{ Cell = 'D8'; Value = '0.596' }
    @{ Cell = 'E8'; Value = '  -0.72%  ' }
    @{ Cell = 'D9'; Value = '0.109' }
    @{ Cell = 'E9'; Value = '  +0.46%  ' }
    @{ Cell = 'D10'; Value = '0.168' }
    @{ Cell = 'E10'; Value = '  +4.99%  ' }
    @{ Cell = 'D11'; Value = '5.69' }
    @{ Cell = 'E11'; Value = '  +0.95%  ' }
    @{ Cell = 'E12'; Value = '  +0.03%  ' }
    @{ Cell = 'D13'; Value = '3.223.86' }
    @{ Cell = 'E13'; Value = '  -0.45%  ' }
    @{ Cell = 'D14'; Value = '27.05' }
    @{ Cell = 'E14'; Value = '  +2.54%  ' }
    @{ Cell = 'D15'; Value = '63.705.51' }
    @{ Cell = 'E15'; Value = '  +0.32%  ' }
    @{ Cell = 'E16'; Value = '  +0.35%  ' }
    @{ Cell = 'D17'; Value = '2.744.25' }
    @{ Cell = 'E17'; Value = '  -0.49%  ' }
    @{ Cell = 'D18'; Value = '12.51' }
    @{ Cell = 'E18'; Value = '  +3.50%  ' }
    @{ Cell = 'E19'; Value = '  -0.78%  ' }
    @{ Cell = 'D20'; Value = '355.97' }
    @{ Cell = 'E20'; Value = '  +0.55%  ' }
    @{ Cell = 'D21'; Value = '6.63' }
    @{ Cell = 'E21'; Value = '  -1.41%  ' }
    @{ Cell = 'D22'; Value = '0.999' }
    @{ Cell = 'E22'; Value = '  +0.12%  ' }
    @{ Cell = 'E23'; Value = '  -1.81%  ' }
    @{ Cell = 'D24'; Value = '64.66' }
    @{ Cell = 'E24'; Value = '  -0.65%  ' }
    @{ Cell = 'E25'; Value = '  +0.34%  ' }
    @{ Cell = 'E26'; Value = '  -0.06%  ' }
    @{ Cell = 'E27'; Value = '  +0.45%  ' }
    @{ Cell = 'D28'; Value = '0.0₃0911' }
    @{ Cell = 'E28'; Value = '  +1.24%  ' }
    @{ Cell = 'D29'; Value = '2.00' }
    @{ Cell = 'E29'; Value = '  +3.90%  ' }
    @{ Cell = 'B30'; Value = 'Fetch.AI' }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet' }
    @{ Cell = 'D30'; Value = '1.35' }
    @{ Cell = 'E30'; Value = '  +12.50%  ' }
    @{ Cell = 'B31'; Value = 'Aptos' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' }
    @{ Cell = 'D31'; Value = '7.19' }
    @{ Cell = 'E31'; Value = '  +2.82%  ' }
    @{ Cell = 'D32'; Value = '166.63' }
    @{ Cell = 'E32'; Value = '  -1.47%  ' }
    @{ Cell = 'D33'; Value = '4.93' }
    @{ Cell = 'E33'; Value = '  +1.72%  ' }
    @{ Cell = 'D34'; Value = '20.15' }
    @{ Cell = 'E34'; Value = '  +0.21%  ' }
    @{ Cell = 'D35'; Value = '1.48' }
    @{ Cell = 'E35'; Value = '  +3.45%  ' }
    @{ Cell = 'D36'; Value = '0.998' }
    @{ Cell = 'E36'; Value = '  -0.02%  ' }
    @{ Cell = 'E37'; Value = '  +1.65%  ' }
    @{ Cell = 'D38'; Value = '0.984' }
    @{ Cell = 'E38'; Value = '  +0.78%  ' }
    @{ Cell = 'D39'; Value = '346.48' }
    @{ Cell = 'E39'; Value = '  +5.96%  ' }
    @{ Cell = 'D40'; Value = '6.33' }
    @{ Cell = 'E40'; Value = '  +2.70%  ' }
    @{ Cell = 'D41'; Value = '4.11' }
    @{ Cell = 'E41'; Value = '  -0.70%  ' }
    @{ Cell = 'D42'; Value = '38.69' }
    @{ Cell = 'E42'; Value = '  -0.82%  ' }
    @{ Cell = 'D43'; Value = '21.82' }
    @{ Cell = 'E43'; Value = '  +2.60%  ' }
    @{ Cell = 'D44'; Value = '21.26' }
    @{ Cell = 'E44'; Value = '  -0.49%  ' }
    @{ Cell = 'E45'; Value = '  +0.75%  ' }
    @{ Cell = 'D46'; Value = '0.631' }
    @{ Cell = 'E46'; Value = '  +1.08%  ' }
    @{ Cell = 'E47'; Value = '  -0.38%  ' }
    @{ Cell = 'E48'; Value = '  -0.31%  ' }
    @{ Cell = 'D49'; Value = '132.38' }
    @{ Cell = 'E49'; Value = '  -1.82%  ' }
    @{ Cell = 'E50'; Value = '  -0.09%  ' }
    @{ Cell = 'D51'; Value = '11.09' }
    @{ Cell = 'E51'; Value = '  +0.38%  ' }
)

foreach ($u in $updates) {
    $text = $u.Value
    if ($u.Cell -match '^D\d+$' -and $text -match '^[+-]?\d+(\.\d+)?$') {
        $ws.Range($u.Cell).Value = "'" + $text
    } else {
        $ws.Range($u.Cell).Value = $text
    }
}
